$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PREPROD claim numbers (replaces the old PREPROD block and extends it)
$preprodClaims = @(
    "0420172010449    ",
    "0420172010457",
    "0420172010451",
    "0420172010448",
    "0420172010450",
    "0420172010452",
    "0420172010453",
    "0420172010454",
    "0420172010455",
    "0420172010458",
    "1120170200983",
    "1220170301481",
    "1220170301482"
)

$row = 6
foreach ($claim in $preprodClaims) {
    $ws.Cells.Item($row, 1).Value = "PREPROD"
    $ws.Cells.Item($row, 2).Value = "'" + $claim
    $row = $row + 1
}

$ws.Range("A17:A18").Select()
